# Insert a new row above row 31 on the "2024" worksheet, shifting the
# existing "Others" group entries (September details in columns R:S,
# overflowing into the August P:Q block, and the trailing "Broadband"
# group header) down by one row, then populate the newly freed row 31
# with the latest September entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows(31).Insert()

$ws.Range("R31").Value = "transfer anyone axis"
$ws.Range("S31").Value = "2024-09-05 16:35:58"
